$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 2350
$ws.Range("I49").Value = 2350
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 7050
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -6914
$ws.Range("N49").ClearContents()
$ws.Range("H112").Value = 1806.3572
$ws.Range("J112").Value = 1806.3572
$ws.Range("L112").Value = 5419.071599999999
$ws.Range("N112").Value = -7635.071599999999
$ws.Range("H137").Value = 1745.2452
$ws.Range("I137").Value = 1454.2162
$ws.Range("J137").Value = 2418.25
$ws.Range("K137").Value = 4362.6486
$ws.Range("L137").Value = 7254.75
$ws.Range("M137").Value = -1812.6486
$ws.Range("N137").Value = -12354.75
$ws.Range("H139").Value = 49999.555
$ws.Range("J139").Value = 49999.555
$ws.Range("L139").Value = 49999.555
$ws.Range("N139").Value = -60279.555
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2900.4827
$ws.Range("I61").Value = 2658.4614
$ws.Range("J61").Value = 4998
$ws.Range("K61").Value = 2658.4614
$ws.Range("L61").Value = 4998
$ws.Range("M61").Value = -2446.4614
$ws.Range("N61").Value = -5422
$ws.Range("H82").Value = 10052
$ws.Range("J82").Value = 10052
$ws.Range("L82").Value = 10052
$ws.Range("N82").Value = -10774
$ws.Range("H85").Value = 10052
$ws.Range("J85").Value = 10052
$ws.Range("L85").Value = 10052
$ws.Range("N85").Value = -12548
$ws.Range("H86").Value = 36665
$ws.Range("J86").Value = 36665
$ws.Range("L86").Value = 36665
$ws.Range("N86").Value = -39037
$ws.Range("H89").Value = 36665
$ws.Range("J89").Value = 36665
$ws.Range("L89").Value = 109995
$ws.Range("N89").Value = -121851
$ws.Range("H102").Value = 1812.5714
$ws.Range("I102").Value = 1901.75
$ws.Range("J102").Value = 1527.2
$ws.Range("K102").Value = 1901.75
$ws.Range("L102").Value = 1527.2
$ws.Range("M102").Value = -279.75
$ws.Range("N102").Value = -4771.2
$ws.Range("H136").Value = 2900.4827
$ws.Range("I136").Value = 2658.4614
$ws.Range("J136").Value = 4998
$ws.Range("K136").Value = 7975.3842
$ws.Range("L136").Value = 14994
$ws.Range("M136").Value = -5425.3842
$ws.Range("N136").Value = -20094
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41248.848
$ws.Range("I31").Value = 1498.579
$ws.Range("K31").Value = 1498.579
$ws.Range("M31").Value = -1203.579
$ws.Range("H34").Value = 41248.848
$ws.Range("I34").Value = 1498.579
$ws.Range("K34").Value = 1498.579
$ws.Range("M34").Value = -1296.579
$ws.Range("H58").Value = 2525.111
$ws.Range("I58").Value = 2564.5
$ws.Range("J58").Value = 2210
$ws.Range("K58").Value = 2564.5
$ws.Range("L58").Value = 2210
$ws.Range("M58").Value = -2361.5
$ws.Range("N58").Value = -2616
$ws.Range("H134").Value = 288604.4
$ws.Range("I134").Value = 2701.6667
$ws.Range("J134").Value = 5005999.5
$ws.Range("K134").Value = 8105.000100000001
$ws.Range("L134").Value = 15017998.5
$ws.Range("M134").Value = -5570.000100000001
$ws.Range("N134").Value = -15023068.5
$ws.Range("H136").Value = 2525.111
$ws.Range("I136").Value = 2564.5
$ws.Range("J136").Value = 2210
$ws.Range("K136").Value = 7693.5
$ws.Range("L136").Value = 6630
$ws.Range("M136").Value = -5143.5
$ws.Range("N136").Value = -11730
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 426014.06
$ws.Range("J132").Value = 628692.3
$ws.Range("L132").Value = 5658230.7
$ws.Range("N132").Value = -5663290.7
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 11990332
$ws.Range("I14").Value = 16005734
$ws.Range("J14").Value = 1550288
$ws.Range("K14").Value = 16005734
$ws.Range("L14").Value = 1550288
$ws.Range("M14").Value = -16005566
$ws.Range("N14").Value = -1550624
$ws.Range("H19").Value = 54499.5
$ws.Range("J19").Value = 99000
$ws.Range("L19").Value = 99000
$ws.Range("N19").Value = -99576
$ws.Range("H20").Value = 9009
$ws.Range("J20").Value = 9009
$ws.Range("L20").Value = 9009
$ws.Range("N20").Value = -9499
$ws.Range("H40").Value = 38499.5
$ws.Range("I40").Value = 10000
$ws.Range("J40").Value = 66999
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 66999
$ws.Range("M40").Value = -9849
$ws.Range("N40").Value = -67301
$ws.Range("H55").Value = 15505.9
$ws.Range("J55").Value = 15838.167
$ws.Range("L55").Value = 15838.167
$ws.Range("N55").Value = -16492.167
$ws.Range("H82").Value = 100000
$ws.Range("J82").Value = 100000
$ws.Range("L82").Value = 100000
$ws.Range("N82").Value = -100766
$ws.Range("H85").Value = 100000
$ws.Range("J85").Value = 100000
$ws.Range("L85").Value = 100000
$ws.Range("N85").Value = -102652
$ws.Range("H86").Value = 86759.664
$ws.Range("J86").Value = 86759.664
$ws.Range("L86").Value = 86759.664
$ws.Range("N86").Value = -89131.664
$ws.Range("H89").Value = 86759.664
$ws.Range("J89").Value = 86759.664
$ws.Range("L89").Value = 260278.992
$ws.Range("N89").Value = -272134.992
$ws.Range("H102").Value = 3565.5652
$ws.Range("I102").Value = 2059.9285
$ws.Range("J102").Value = 5907.6665
$ws.Range("K102").Value = 2059.9285
$ws.Range("L102").Value = 5907.6665
$ws.Range("M102").Value = -437.9285
$ws.Range("N102").Value = -9151.666499999999
$ws.Range("H122").Value = 2632.8696
$ws.Range("I122").Value = 2034.875
$ws.Range("J122").Value = 3999.7144
$ws.Range("K122").Value = 6104.625
$ws.Range("L122").Value = 11999.1432
$ws.Range("M122").Value = -3654.625
$ws.Range("N122").Value = -16899.1432
$ws.Range("H126").Value = 2855
$ws.Range("J126").Value = 3200.5
$ws.Range("L126").Value = 9601.5
$ws.Range("N126").Value = -14541.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 13300
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H46").Value = 2069
$ws.Range("I46").Value = 2206.75
$ws.Range("K46").Value = 2206.75
$ws.Range("M46").Value = -2018.75
$ws.Range("H122").Value = 5503.7896
$ws.Range("I122").Value = 5155.2144
$ws.Range("J122").Value = 6479.8
$ws.Range("K122").Value = 15465.6432
$ws.Range("L122").Value = 19439.4
$ws.Range("M122").Value = -13015.6432
$ws.Range("N122").Value = -24339.4
$ws.Range("H127").Value = 83173.25
$ws.Range("J127").Value = 83173.25
$ws.Range("L127").Value = 83173.25
$ws.Range("N127").Value = -93093.25
$ws.Range("H132").Value = 3732.5386
$ws.Range("I132").Value = 3543.5833
$ws.Range("K132").Value = 10630.7499
$ws.Range("M132").Value = -8100.749899999999
$ws.Range("H135").Value = 69507.56
$ws.Range("J135").Value = 69507.56
$ws.Range("L135").Value = 69507.56
$ws.Range("N135").Value = -79647.56
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12018.75
$ws.Range("I41").Value = 11578
$ws.Range("J41").Value = 12165.667
$ws.Range("K41").Value = 11578
$ws.Range("L41").Value = 12165.667
$ws.Range("M41").Value = -11188
$ws.Range("N41").Value = -12945.667
$ws.Range("H56").Value = 50000
$ws.Range("J56").Value = 50000
$ws.Range("L56").Value = 50000
$ws.Range("N56").Value = -51428
$ws.Range("H82").Value = 52836.855
$ws.Range("J82").Value = 60571.6
$ws.Range("L82").Value = 60571.6
$ws.Range("N82").Value = -61337.6
$ws.Range("H85").Value = 52836.855
$ws.Range("J85").Value = 60571.6
$ws.Range("L85").Value = 60571.6
$ws.Range("N85").Value = -63223.6
$ws.Range("H110").Value = 100000
$ws.Range("J110").Value = 100000
$ws.Range("L110").Value = 100000
$ws.Range("N110").Value = -108180
$ws.Range("H126").Value = 1895.7858
$ws.Range("I126").Value = 1899.6364
$ws.Range("K126").Value = 5698.9092
$ws.Range("M126").Value = -3228.9092

Write-Host "Applied all cell updates"